$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-format all the hour/minute cells (rows 2-15, columns C-Q) from the old
# "h:mm;@" display format to the correct "[h]:mm" elapsed-time format used
# for time-budget calculations.
$ws.Range("C2:Q15").NumberFormat = "[h]:mm"

# Two new rows were appended below the table; give them the same time
# format. A18 additionally carries a small Consolas font.
$ws.Range("A18").Font.Name = "Consolas"
$ws.Range("A18").Font.Family = 3
$ws.Range("A18").Font.Size = 10
$ws.Range("A18").Font.Color = 2499618
$ws.Range("A18").NumberFormat = "[h]:mm"

$ws.Range("A19").NumberFormat = "[h]:mm"

$ws.Range("D17").Select()
